$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old "TOTAL" row (old row 15), pushing it down to row 17.
# This mirrors what the author did in Excel: the grid grew from 5 product rows (12-14 + total)
# to 5 product rows (12-16) + a total row (17).
$ws.Rows("15:16").Insert()

# --- Row 12: Pain complet / A la pièce ---
$ws.Range("A12").Value = "Pain complet"
$ws.Range("B12").Value = "Pain complet"
$ws.Range("C12").Value = "A la pièce"
$ws.Range("D12").Value = 4
$ws.Range("E12").Value = 1
$ws.Range("F12").Formula = "=D12*E12"

# --- Row 13: Pomme de terre / Au poids ---
$ws.Range("A13").Value = "Pomme de terre"
$ws.Range("B13").Value = "Pomme de terre"
$ws.Range("C13").Value = "Au poids"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 3
$ws.Range("F13").Formula = "=D13*E13"

# --- Row 14: Radis / A la pièce ---
$ws.Range("A14").Value = "Radis"
$ws.Range("B14").Value = "Radis"
$ws.Range("C14").Value = "A la pièce"
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 2
$ws.Range("F14").Formula = "=D14*E14"

# --- Row 15 (new): Salade / A la pièce ---
$ws.Range("A15").Value = "Salade"
$ws.Range("B15").Value = "Salade"
$ws.Range("C15").Value = "A la pièce"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 2
$ws.Range("F15").Formula = "=D15*E15"
$ws.Range("A15:B15").Merge()

# --- Row 16 (new): Tomates grappe / Au poids ---
$ws.Range("A16").Value = "Tomates grappe"
$ws.Range("B16").Value = "Tomates grappe"
$ws.Range("C16").Value = "Au poids"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 2
$ws.Range("F16").Formula = "=D16*E16"
$ws.Range("A16:B16").Merge()

# --- Row 17: TOTAL row (was row 15, shifted down by the insert above) ---
$ws.Range("E17").Value = "TOTAL : "
$ws.Range("F17").Formula = "=SUBTOTAL(9,F12:F16)"
